# Update NATMI ligand/receptor TPM-derived expression & specificity values
# (Pdgfb-Lrp1 sheet) to reflect the new TPM normalization.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 32.544082
$ws.Range("H2").Value = 97.632246
$ws.Range("I2").Value = 0.621589875979724
$ws.Range("J2").Value = 0.6366365948489335
$ws.Range("M2").Value = 2.906846333333333
$ws.Range("N2").Value = 8.720538999999999
$ws.Range("O2").Value = 0.005520525738044089
$ws.Range("P2").Value = 0.005624540846623205
$ws.Range("Q2").Value = 94.60064543339931
$ws.Range("R2").Value = 851.4058089005938
$ws.Range("S2").Value = 0.0034315029088537
$ws.Range("T2").Value = 0.003580788532182935
$ws.Range("G3").Value = 32.544082
$ws.Range("H3").Value = 97.632246
$ws.Range("I3").Value = 0.621589875979724
$ws.Range("J3").Value = 0.6366365948489335
$ws.Range("O3").Value = 0.3528665483720876
$ws.Range("P3").Value = 0.3595150912979765
$ws.Range("Q3").Value = 6046.779747409034
$ws.Range("R3").Value = 54421.01772668131
$ws.Range("S3").Value = 0.2193382740399993
$ws.Range("T3").Value = 0.2288804635207472
$ws.Range("G4").Value = 32.544082
$ws.Range("H4").Value = 97.632246
$ws.Range("I4").Value = 0.621589875979724
$ws.Range("J4").Value = 0.6366365948489335
$ws.Range("M4").Value = 137.0717086666666
$ws.Range("N4").Value = 411.2151259999999
$ws.Range("O4").Value = 0.2603191943704447
$ws.Range("P4").Value = 0.2652240042658267
$ws.Range("Q4").Value = 4460.872926728109
$ws.Range("R4").Value = 40147.85634055299
$ws.Range("S4").Value = 0.1618117757438664
$ws.Range("T4").Value = 0.1688513069479949
$ws.Range("G5").Value = 32.544082
$ws.Range("H5").Value = 97.632246
$ws.Range("I5").Value = 0.621589875979724
$ws.Range("J5").Value = 0.6366365948489335
$ws.Range("M5").Value = 29.2127365
$ws.Range("N5").Value = 58.425473
$ws.Range("O5").Value = 0.05547925319534149
$ws.Range("P5").Value = 0.03768304451958546
$ws.Range("Q5").Value = 950.7016921003928
$ws.Range("R5").Value = 5704.210152602357
$ws.Range("S5").Value = 0.03448534211314002
$ws.Range("T5").Value = 0.02399040514648965
$ws.Range("G6").Value = 32.544082
$ws.Range("H6").Value = 97.632246
$ws.Range("I6").Value = 0.621589875979724
$ws.Range("J6").Value = 0.6366365948489335
$ws.Range("M6").Value = 171.5584106666666
$ws.Range("N6").Value = 514.6752319999999
$ws.Range("O6").Value = 0.3258144783240821
$ws.Range("P6").Value = 0.331953319069988
$ws.Range("Q6").Value = 5583.210984525674
$ws.Range("R6").Value = 50248.89886073106
$ws.Range("S6").Value = 0.2025229811738647
$ws.Range("T6").Value = 0.2113336307015187
$ws.Range("I7").Value = 0.004665102012661462
$ws.Range("J7").Value = 0.004778029332093849
$ws.Range("M7").Value = 2.906846333333333
$ws.Range("N7").Value = 8.720538999999999
$ws.Range("O7").Value = 0.005520525738044089
$ws.Range("P7").Value = 0.005624540846623205
$ws.Range("Q7").Value = 0.7099884963776665
$ws.Range("R7").Value = 6.389896467398999
$ws.Range("S7").Value = 0.00002575381573149888
$ws.Range("T7").Value = 0.00002687422114472565
$ws.Range("I8").Value = 0.004665102012661462
$ws.Range("J8").Value = 0.004778029332093849
$ws.Range("O8").Value = 0.3528665483720876
$ws.Range("P8").Value = 0.3595150912979765
$ws.Range("S8").Value = 0.001646158445011529
$ws.Range("T8").Value = 0.00171777365155213
$ws.Range("I9").Value = 0.004665102012661462
$ws.Range("J9").Value = 0.004778029332093849
$ws.Range("M9").Value = 137.0717086666666
$ws.Range("N9").Value = 411.2151259999999
$ws.Range("O9").Value = 0.2603191943704447
$ws.Range("P9").Value = 0.2652240042658267
$ws.Range("Q9").Value = 33.47935362670732
$ws.Range("R9").Value = 301.314182640366
$ws.Range("S9").Value = 0.001214415597591972
$ws.Range("T9").Value = 0.001267248071957504
$ws.Range("I10").Value = 0.004665102012661462
$ws.Range("J10").Value = 0.004778029332093849
$ws.Range("M10").Value = 29.2127365
$ws.Range("N10").Value = 58.425473
$ws.Range("O10").Value = 0.05547925319534149
$ws.Range("P10").Value = 0.03768304451958546
$ws.Range("Q10").Value = 7.135123251915499
$ws.Range("R10").Value = 42.810739511493
$ws.Range("S10").Value = 0.0002588163757425424
$ws.Range("T10").Value = 0.0001800506920371777
$ws.Range("I11").Value = 0.004665102012661462
$ws.Range("J11").Value = 0.004778029332093849
$ws.Range("M11").Value = 171.5584106666666
$ws.Range("N11").Value = 514.6752319999999
$ws.Range("O11").Value = 0.3258144783240821
$ws.Range("P11").Value = 0.331953319069988
$ws.Range("Q11").Value = 41.90262713010133
$ws.Range("R11").Value = 377.1236441709119
$ws.Range("S11").Value = 0.00151995777858392
$ws.Range("T11").Value = 0.001586082695402311
$ws.Range("G12").Value = 6.619872666666667
$ws.Range("H12").Value = 19.859618
$ws.Range("I12").Value = 0.126439142756428
$ws.Range("J12").Value = 0.1294998332673878
$ws.Range("M12").Value = 2.906846333333333
$ws.Range("N12").Value = 8.720538999999999
$ws.Range("O12").Value = 0.005520525738044089
$ws.Range("P12").Value = 0.005624540846623205
$ws.Range("Q12").Value = 19.24295258823355
$ws.Range("R12").Value = 173.186573294102
$ws.Range("S12").Value = 0.0006980105418830916
$ws.Range("T12").Value = 0.0007283771018433173
$ws.Range("G13").Value = 6.619872666666667
$ws.Range("H13").Value = 19.859618
$ws.Range("I13").Value = 0.126439142756428
$ws.Range("J13").Value = 0.1294998332673878
$ws.Range("O13").Value = 0.3528665483720876
$ws.Range("P13").Value = 0.3595150912979765
$ws.Range("Q13").Value = 1229.990508603888
$ws.Range("R13").Value = 11069.91457743499
$ws.Range("S13").Value = 0.0446161438835864
$ws.Range("T13").Value = 0.04655714438019766
$ws.Range("G14").Value = 6.619872666666667
$ws.Range("H14").Value = 19.859618
$ws.Range("I14").Value = 0.126439142756428
$ws.Range("J14").Value = 0.1294998332673878
$ws.Range("M14").Value = 137.0717086666666
$ws.Range("N14").Value = 411.2151259999999
$ws.Range("O14").Value = 0.2603191943704447
$ws.Range("P14").Value = 0.2652240042658267
$ws.Range("Q14").Value = 907.397257575763
$ws.Range("R14").Value = 8166.575318181867
$ws.Range("S14").Value = 0.03291453577924298
$ws.Range("T14").Value = 0.0343464643309335
$ws.Range("G15").Value = 6.619872666666667
$ws.Range("H15").Value = 19.859618
$ws.Range("I15").Value = 0.126439142756428
$ws.Range("J15").Value = 0.1294998332673878
$ws.Range("M15").Value = 29.2127365
$ws.Range("N15").Value = 58.425473
$ws.Range("O15").Value = 0.05547925319534149
$ws.Range("P15").Value = 0.03768304451958546
$ws.Range("Q15").Value = 193.3845958748857
$ws.Range("R15").Value = 1160.307575249314
$ws.Range("S15").Value = 0.007014749214785797
$ws.Range("T15").Value = 0.004879947982293868
$ws.Range("G16").Value = 6.619872666666667
$ws.Range("H16").Value = 19.859618
$ws.Range("I16").Value = 0.126439142756428
$ws.Range("J16").Value = 0.1294998332673878
$ws.Range("M16").Value = 171.5584106666666
$ws.Range("N16").Value = 514.6752319999999
$ws.Range("O16").Value = 0.3258144783240821
$ws.Range("P16").Value = 0.331953319069988
$ws.Range("Q16").Value = 1135.694833509042
$ws.Range("R16").Value = 10221.25350158138
$ws.Range("S16").Value = 0.04119570333692973
$ws.Range("T16").Value = 0.04298789947211943
$ws.Range("G17").Value = 3.71227
$ws.Range("H17").Value = 7.42454
$ws.Range("I17").Value = 0.07090411857072049
$ws.Range("J17").Value = 0.04841365488938666
$ws.Range("M17").Value = 2.906846333333333
$ws.Range("N17").Value = 8.720538999999999
$ws.Range("O17").Value = 0.005520525738044089
$ws.Range("P17").Value = 0.005624540846623205
$ws.Range("Q17").Value = 10.79099843784333
$ws.Range("R17").Value = 64.74599062706
$ws.Range("S17").Value = 0.0003914280115029923
$ws.Range("T17").Value = 0.0002723045794596745
$ws.Range("G18").Value = 3.71227
$ws.Range("H18").Value = 7.42454
$ws.Range("I18").Value = 0.07090411857072049
$ws.Range("J18").Value = 0.04841365488938666
$ws.Range("O18").Value = 0.3528665483720876
$ws.Range("P18").Value = 0.3595150912979765
$ws.Range("Q18").Value = 689.7499537063034
$ws.Range("R18").Value = 4138.499722237821
$ws.Range("S18").Value = 0.02501969158541538
$ws.Range("T18").Value = 0.01740543955762658
$ws.Range("G19").Value = 3.71227
$ws.Range("H19").Value = 7.42454
$ws.Range("I19").Value = 0.07090411857072049
$ws.Range("J19").Value = 0.04841365488938666
$ws.Range("M19").Value = 137.0717086666666
$ws.Range("N19").Value = 411.2151259999999
$ws.Range("O19").Value = 0.2603191943704447
$ws.Range("P19").Value = 0.2652240042658267
$ws.Range("Q19").Value = 508.8471919320066
$ws.Range("R19").Value = 3053.083151592039
$ws.Range("S19").Value = 0.01845770302387644
$ws.Range("T19").Value = 0.01284046341090695
$ws.Range("G20").Value = 3.71227
$ws.Range("H20").Value = 7.42454
$ws.Range("I20").Value = 0.07090411857072049
$ws.Range("J20").Value = 0.04841365488938666
$ws.Range("M20").Value = 29.2127365
$ws.Range("N20").Value = 58.425473
$ws.Range("O20").Value = 0.05547925319534149
$ws.Range("P20").Value = 0.03768304451958546
$ws.Range("Q20").Value = 108.445565326855
$ws.Range("R20").Value = 433.78226130742
$ws.Range("S20").Value = 0.003933707546777516
$ws.Range("T20").Value = 0.001824373912552604
$ws.Range("G21").Value = 3.71227
$ws.Range("H21").Value = 7.42454
$ws.Range("I21").Value = 0.07090411857072049
$ws.Range("J21").Value = 0.04841365488938666
$ws.Range("M21").Value = 171.5584106666666
$ws.Range("N21").Value = 514.6752319999999
$ws.Range("O21").Value = 0.3258144783240821
$ws.Range("P21").Value = 0.331953319069988
$ws.Range("Q21").Value = 636.8711411655466
$ws.Range("R21").Value = 3821.22684699328
$ws.Range("S21").Value = 0.02310158840314816
$ws.Range("T21").Value = 0.01607107342884086
$ws.Range("G22").Value = 9.235725333333333
$ws.Range("H22").Value = 27.707176
$ws.Range("I22").Value = 0.176401760680466
$ws.Range("J22").Value = 0.1806718876621981
$ws.Range("M22").Value = 2.906846333333333
$ws.Range("N22").Value = 8.720538999999999
$ws.Range("O22").Value = 0.005520525738044089
$ws.Range("P22").Value = 0.005624540846623205
$ws.Range("Q22").Value = 26.84683432087377
$ws.Range("R22").Value = 241.6215088878639
$ws.Range("S22").Value = 0.0009738304600728064
$ws.Range("T22").Value = 0.001016196411992553
$ws.Range("G23").Value = 9.235725333333333
$ws.Range("H23").Value = 27.707176
$ws.Range("I23").Value = 0.176401760680466
$ws.Range("J23").Value = 0.1806718876621981
$ws.Range("O23").Value = 0.3528665483720876
$ws.Range("P23").Value = 0.3595150912979765
$ws.Range("Q23").Value = 1716.023112842223
$ws.Range("R23").Value = 15444.20801558001
$ws.Range("S23").Value = 0.0622462804180751
$ws.Range("T23").Value = 0.06495427018785292
$ws.Range("G24").Value = 9.235725333333333
$ws.Range("H24").Value = 27.707176
$ws.Range("I24").Value = 0.176401760680466
$ws.Range("J24").Value = 0.1806718876621981
$ws.Range("M24").Value = 137.0717086666666
$ws.Range("N24").Value = 411.2151259999999
$ws.Range("O24").Value = 0.2603191943704447
$ws.Range("P24").Value = 0.2652240042658267
$ws.Range("Q24").Value = 1265.956652216019
$ws.Range("R24").Value = 11393.60986994417
$ws.Range("S24").Value = 0.0459207642258669
$ws.Range("T24").Value = 0.0479185215040338
$ws.Range("G25").Value = 9.235725333333333
$ws.Range("H25").Value = 27.707176
$ws.Range("I25").Value = 0.176401760680466
$ws.Range("J25").Value = 0.1806718876621981
$ws.Range("M25").Value = 29.2127365
$ws.Range("N25").Value = 58.425473
$ws.Range("O25").Value = 0.05547925319534149
$ws.Range("P25").Value = 0.03768304451958546
$ws.Range("Q25").Value = 269.8008105490413
$ws.Range("R25").Value = 1618.804863294248
$ws.Range("S25").Value = 0.009786637944895609
$ws.Range("T25").Value = 0.006808266786212155
$ws.Range("G26").Value = 9.235725333333333
$ws.Range("H26").Value = 27.707176
$ws.Range("I26").Value = 0.176401760680466
$ws.Range("J26").Value = 0.1806718876621981
$ws.Range("M26").Value = 171.5584106666666
$ws.Range("N26").Value = 514.6752319999999
$ws.Range("O26").Value = 0.3258144783240821
$ws.Range("P26").Value = 0.331953319069988
$ws.Range("Q26").Value = 1584.466359540537
$ws.Range("R26").Value = 14260.19723586483
$ws.Range("S26").Value = 0.05747424763155562
$ws.Range("T26").Value = 0.05997463277210669
